$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =====================================================================
# Rows 50/51: the two matches on 2023-08-11 swap places (id in column A
# stays tied to the row position; all other fields move with the match).
# =====================================================================
# Row 50
$ws.Range("B50").Value = 6732795
$ws.Range("C50").Value = "Lithuania A Lyga"
$ws.Range("D50").Value = "Lithuania A Lyga"
$ws.Range("F50").Value = "Suduva Marijampole"
$ws.Range("G50").Value = "Banga Gargzdai"
$ws.Range("H50").Value = 1
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = "H"
$ws.Range("K50").Value = 2.15
$ws.Range("L50").Value = 3.2
$ws.Range("M50").Value = 3
$ws.Range("N50").Value = 2.3
$ws.Range("O50").Value = 3.2
$ws.Range("P50").Value = 2.7
$ws.Range("Q50").Value = -0.25
$ws.Range("R50").Value = 2.05
$ws.Range("S50").Value = 1.75
$ws.Range("T50").Value = 2.25
$ws.Range("U50").Value = 1.9
$ws.Range("V50").Value = 1.9
$ws.Range("W50").Value = 1.3
$ws.Range("X50").Value = -1
$ws.Range("Y50").Value = -1
$ws.Range("Z50").Value = 1.05
$ws.Range("AA50").Value = -1
$ws.Range("AB50").Value = -1
$ws.Range("AC50").Value = 0.8999999999999999

# Row 51
$ws.Range("B51").Value = 6732794
$ws.Range("C51").Value = "Lithuania A Lyga"
$ws.Range("D51").Value = "Lithuania A Lyga"
$ws.Range("F51").Value = "FK Siauliai"
$ws.Range("G51").Value = "FK Dziugas Telsiai"
$ws.Range("H51").Value = 3
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = "H"
$ws.Range("K51").Value = 1.25
$ws.Range("L51").Value = 5
$ws.Range("M51").Value = 9
$ws.Range("N51").Value = 1.25
$ws.Range("O51").Value = 5.25
$ws.Range("P51").Value = 9
$ws.Range("Q51").Value = -1.75
$ws.Range("R51").Value = 2
$ws.Range("S51").Value = 1.8
$ws.Range("T51").Value = 3
$ws.Range("U51").Value = 1.975
$ws.Range("V51").Value = 1.825
$ws.Range("W51").Value = 0.25
$ws.Range("X51").Value = -1
$ws.Range("Y51").Value = -1
$ws.Range("Z51").Value = 1
$ws.Range("AA51").Value = -1
$ws.Range("AB51").Value = 0
$ws.Range("AC51").Value = 0

# =====================================================================
# Rows 100-104: the five matches on 2023-11-12 are re-ordered (id in
# column A stays tied to the row position).
# =====================================================================
# Row 100
$ws.Range("B100").Value = 6732836
$ws.Range("C100").Value = "Lithuania A Lyga"
$ws.Range("D100").Value = "Lithuania A Lyga"
$ws.Range("F100").Value = "FK Siauliai"
$ws.Range("G100").Value = "Banga Gargzdai"
$ws.Range("H100").Value = 3
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = "H"
$ws.Range("K100").Value = 1.222
$ws.Range("L100").Value = 5.5
$ws.Range("M100").Value = 9
$ws.Range("N100").Value = 1.363
$ws.Range("O100").Value = 4.5
$ws.Range("P100").Value = 7
$ws.Range("Q100").Value = -1.25
$ws.Range("R100").Value = 1.9
$ws.Range("S100").Value = 1.9
$ws.Range("T100").Value = 2.5
$ws.Range("U100").Value = 1.975
$ws.Range("V100").Value = 1.825
$ws.Range("W100").Value = 0.363
$ws.Range("X100").Value = -1
$ws.Range("Y100").Value = -1
$ws.Range("Z100").Value = 0.8999999999999999
$ws.Range("AA100").Value = -1
$ws.Range("AB100").Value = 0.9750000000000001
$ws.Range("AC100").Value = -1

# Row 101
$ws.Range("B101").Value = 6732727
$ws.Range("C101").Value = "Lithuania A Lyga"
$ws.Range("D101").Value = "Lithuania A Lyga"
$ws.Range("F101").Value = "FK Zalgiris Vilnius"
$ws.Range("G101").Value = "FK Dainava Alytus"
$ws.Range("H101").Value = 1
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = "H"
$ws.Range("K101").Value = 1.285
$ws.Range("L101").Value = 5.5
$ws.Range("M101").Value = 6.5
$ws.Range("N101").Value = 1.3
$ws.Range("O101").Value = 5.5
$ws.Range("P101").Value = 6
$ws.Range("Q101").Value = -1.5
$ws.Range("R101").Value = 1.9
$ws.Range("S101").Value = 1.9
$ws.Range("T101").Value = 2.75
$ws.Range("U101").Value = 1.8
$ws.Range("V101").Value = 2
$ws.Range("W101").Value = 0.3
$ws.Range("X101").Value = -1
$ws.Range("Y101").Value = -1
$ws.Range("Z101").Value = -1
$ws.Range("AA101").Value = 0.8999999999999999
$ws.Range("AB101").Value = -1
$ws.Range("AC101").Value = 1

# Row 102
$ws.Range("B102").Value = 6732837
$ws.Range("C102").Value = "Lithuania A Lyga"
$ws.Range("D102").Value = "Lithuania A Lyga"
$ws.Range("F102").Value = "Suduva Marijampole"
$ws.Range("G102").Value = "FK Riteriai"
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 3
$ws.Range("J102").Value = "A"
$ws.Range("K102").Value = 3.6
$ws.Range("L102").Value = 3.6
$ws.Range("M102").Value = 1.8
$ws.Range("N102").Value = 3
$ws.Range("O102").Value = 3.6
$ws.Range("P102").Value = 2
$ws.Range("Q102").Value = 0.25
$ws.Range("R102").Value = 2
$ws.Range("S102").Value = 1.8
$ws.Range("T102").Value = 2.5
$ws.Range("U102").Value = 1.975
$ws.Range("V102").Value = 1.825
$ws.Range("W102").Value = -1
$ws.Range("X102").Value = -1
$ws.Range("Y102").Value = 1
$ws.Range("Z102").Value = -1
$ws.Range("AA102").Value = 0.8
$ws.Range("AB102").Value = 0.9750000000000001
$ws.Range("AC102").Value = -1

# Row 103
$ws.Range("B103").Value = 7465686
$ws.Range("C103").Value = "Lithuania A Lyga"
$ws.Range("D103").Value = "Lithuania A Lyga"
$ws.Range("F103").Value = "FK Kauno Zalgiris"
$ws.Range("G103").Value = "Hegelmann Litauen"
$ws.Range("H103").Value = 4
$ws.Range("I103").Value = 2
$ws.Range("J103").Value = "H"
$ws.Range("K103").Value = 2.3
$ws.Range("L103").Value = 4
$ws.Range("M103").Value = 2.3
$ws.Range("N103").Value = 2.55
$ws.Range("O103").Value = 4
$ws.Range("P103").Value = 2.2
$ws.Range("Q103").Value = 0.25
$ws.Range("R103").Value = 1.8
$ws.Range("S103").Value = 2
$ws.Range("T103").Value = 2.75
$ws.Range("U103").Value = 1.85
$ws.Range("V103").Value = 1.95
$ws.Range("W103").Value = 1.55
$ws.Range("X103").Value = -1
$ws.Range("Y103").Value = -1
$ws.Range("Z103").Value = 0.8
$ws.Range("AA103").Value = -1
$ws.Range("AB103").Value = 0.8500000000000001
$ws.Range("AC103").Value = -1

# Row 104
$ws.Range("B104").Value = 6732834
$ws.Range("C104").Value = "Lithuania A Lyga"
$ws.Range("D104").Value = "Lithuania A Lyga"
$ws.Range("F104").Value = "Panevezys"
$ws.Range("G104").Value = "FK Dziugas Telsiai"
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = "D"
$ws.Range("K104").Value = 1.25
$ws.Range("L104").Value = 5.5
$ws.Range("M104").Value = 7.5
$ws.Range("N104").Value = 1.45
$ws.Range("O104").Value = 4.5
$ws.Range("P104").Value = 5
$ws.Range("Q104").Value = -1
$ws.Range("R104").Value = 1.775
$ws.Range("S104").Value = 2.025
$ws.Range("T104").Value = 2.5
$ws.Range("U104").Value = 1.875
$ws.Range("V104").Value = 1.925
$ws.Range("W104").Value = -1
$ws.Range("X104").Value = 3.5
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = -1
$ws.Range("AA104").Value = 1.025
$ws.Range("AB104").Value = -1
$ws.Range("AC104").Value = 0.925

# =====================================================================
# Rows 105-108: four new matches appended at the bottom of the table.
# Column A/E formatting (bold+border / date format) is copied from an
# existing data row so the new rows match the sheet's look.
# =====================================================================
# Row 105
$ws.Range("A50").Copy()
$ws.Range("A105").PasteSpecial(-4122)
$ws.Range("E50").Copy()
$ws.Range("E105").PasteSpecial(-4122)
$ws.Range("A105").Value = 103
$ws.Range("B105").Value = 7862033
$ws.Range("C105").Value = "Lithuania A Lyga"
$ws.Range("D105").Value = "Lithuania A Lyga"
$ws.Range("E105").Value = 45352.58333333334
$ws.Range("F105").Value = "Suduva Marijampole"
$ws.Range("G105").Value = "FK Siauliai"
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = "D"
$ws.Range("K105").Value = 4.333
$ws.Range("L105").Value = 3.3
$ws.Range("M105").Value = 1.727
$ws.Range("N105").Value = 4
$ws.Range("O105").Value = 3.2
$ws.Range("P105").Value = 1.833
$ws.Range("Q105").Value = 0.5
$ws.Range("R105").Value = 1.925
$ws.Range("S105").Value = 1.875
$ws.Range("T105").Value = 2
$ws.Range("U105").Value = 1.775
$ws.Range("V105").Value = 2.025
$ws.Range("W105").Value = -1
$ws.Range("X105").Value = 2.2
$ws.Range("Y105").Value = -1
$ws.Range("Z105").Value = 0.925
$ws.Range("AA105").Value = -1
$ws.Range("AB105").Value = -1
$ws.Range("AC105").Value = 1.025

# Row 106
$ws.Range("A50").Copy()
$ws.Range("A106").PasteSpecial(-4122)
$ws.Range("E50").Copy()
$ws.Range("E106").PasteSpecial(-4122)
$ws.Range("A106").Value = 104
$ws.Range("B106").Value = 7862902
$ws.Range("C106").Value = "Lithuania A Lyga"
$ws.Range("D106").Value = "Lithuania A Lyga"
$ws.Range("E106").Value = 45353.33333333334
$ws.Range("F106").Value = "FK Dziugas Telsiai"
$ws.Range("G106").Value = "Panevezys"
$ws.Range("K106").Value = 7.5
$ws.Range("L106").Value = 4.333
$ws.Range("M106").Value = 1.333
$ws.Range("N106").Value = 7
$ws.Range("O106").Value = 4.333
$ws.Range("P106").Value = 1.363
$ws.Range("Q106").Value = 1.25
$ws.Range("R106").Value = 1.95
$ws.Range("S106").Value = 1.85
$ws.Range("T106").Value = 2.25
$ws.Range("U106").Value = 1.85
$ws.Range("V106").Value = 1.95
$ws.Range("W106").Value = 0
$ws.Range("X106").Value = 0
$ws.Range("Y106").Value = 0
$ws.Range("Z106").Value = 0
$ws.Range("AA106").Value = 0

# Row 107
$ws.Range("A50").Copy()
$ws.Range("A107").PasteSpecial(-4122)
$ws.Range("E50").Copy()
$ws.Range("E107").PasteSpecial(-4122)
$ws.Range("A107").Value = 105
$ws.Range("B107").Value = 7862903
$ws.Range("C107").Value = "Lithuania A Lyga"
$ws.Range("D107").Value = "Lithuania A Lyga"
$ws.Range("E107").Value = 45353.41666666666
$ws.Range("F107").Value = "Banga Gargzdai"
$ws.Range("G107").Value = "FK Kauno Zalgiris"
$ws.Range("K107").Value = 5
$ws.Range("L107").Value = 3.6
$ws.Range("M107").Value = 1.571
$ws.Range("N107").Value = 5.75
$ws.Range("O107").Value = 3.75
$ws.Range("P107").Value = 1.5
$ws.Range("Q107").Value = 1
$ws.Range("R107").Value = 1.95
$ws.Range("S107").Value = 1.85
$ws.Range("T107").Value = 2.5
$ws.Range("U107").Value = 2.025
$ws.Range("V107").Value = 1.775
$ws.Range("W107").Value = 0
$ws.Range("X107").Value = 0
$ws.Range("Y107").Value = 0
$ws.Range("Z107").Value = 0
$ws.Range("AA107").Value = 0

# Row 108
$ws.Range("A50").Copy()
$ws.Range("A108").PasteSpecial(-4122)
$ws.Range("E50").Copy()
$ws.Range("E108").PasteSpecial(-4122)
$ws.Range("A108").Value = 106
$ws.Range("B108").Value = 7862034
$ws.Range("C108").Value = "Lithuania A Lyga"
$ws.Range("D108").Value = "Lithuania A Lyga"
$ws.Range("E108").Value = 45354.33333333334
$ws.Range("F108").Value = "FK Dainava Alytus"
$ws.Range("G108").Value = "FK Zalgiris Vilnius"
$ws.Range("K108").Value = 6.5
$ws.Range("L108").Value = 4
$ws.Range("M108").Value = 1.4
$ws.Range("N108").Value = 7.5
$ws.Range("O108").Value = 4.5
$ws.Range("P108").Value = 1.3
$ws.Range("Q108").Value = 1.25
$ws.Range("R108").Value = 2
$ws.Range("S108").Value = 1.8
$ws.Range("T108").Value = 2.5
$ws.Range("U108").Value = 1.8
$ws.Range("V108").Value = 2
$ws.Range("W108").Value = 0
$ws.Range("X108").Value = 0
$ws.Range("Y108").Value = 0
$ws.Range("Z108").Value = 0
$ws.Range("AA108").Value = 0

$ws.Application.CutCopyMode = $false
